$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-12 for columns I (I0) and J (IF)
$data = @(
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(1, 3),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(6, 7),
    @(4, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
